$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 186.25
$ws.Range("I2").Value = 148.33333
$ws.Range("J2").Value = 300
$ws.Range("K2").Value = 148.33333
$ws.Range("L2").Value = 300
$ws.Range("M2").Value = -35.33332999999999
$ws.Range("N2").Value = -526

$ws.Range("H28").Value = 51162.9
$ws.Range("I28").Value = 68112.07000000001
$ws.Range("K28").Value = 68112.07000000001
$ws.Range("M28").Value = -67627.07000000001

$ws.Range("H107").Value = 40695.76
$ws.Range("I107").Value = 40695.76
$ws.Range("K107").Value = 40695.76
$ws.Range("M107").Value = -38775.76

$ws.Range("H125").Value = 11115654

$ws.Range("H132").Value = 1329.1111
$ws.Range("I132").Value = 1405.2245
$ws.Range("K132").Value = 4215.6735
$ws.Range("M132").Value = -1685.6735

$ws.Range("H137").Value = 3701.5
$ws.Range("I137").Value = 3347.5454
$ws.Range("K137").Value = 10042.6362
$ws.Range("M137").Value = -7492.636200000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5631.647
$ws.Range("I32").Value = 5777.2153
$ws.Range("K32").Value = 5777.2153
$ws.Range("M32").Value = -5490.2153

$ws.Range("H74").Value = 1639.8096
$ws.Range("I74").Value = 1639.8096
$ws.Range("K74").Value = 1639.8096
$ws.Range("M74").Value = -765.8096

$ws.Range("H77").Value = 1639.8096
$ws.Range("I77").Value = 1639.8096
$ws.Range("K77").Value = 8199.048000000001
$ws.Range("M77").Value = -3831.048000000001

$ws.Range("H110").Value = 103349.41
$ws.Range("I110").Value = 107698.32
$ws.Range("J110").Value = 1150
$ws.Range("K110").Value = 107698.32
$ws.Range("L110").Value = 1150
$ws.Range("M110").Value = -105653.32
$ws.Range("N110").Value = -5240

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41669436
$ws.Range("I20").Value = 55557660
$ws.Range("K20").Value = 55557660
$ws.Range("M20").Value = -55557413

$ws.Range("H86").Value = 1001699.06
$ws.Range("I86").Value = 1309329.5
$ws.Range("K86").Value = 1309329.5
$ws.Range("M86").Value = -1308206.5

$ws.Range("H89").Value = 1001699.06
$ws.Range("I89").Value = 1309329.5
$ws.Range("K89").Value = 6546647.5
$ws.Range("M89").Value = -6541031.5

$ws.Range("H94").Value = 235.2
$ws.Range("I94").Value = 218.75
$ws.Range("K94").Value = 218.75
$ws.Range("M94").Value = 232.25

$ws.Range("H134").Value = 37706.2
$ws.Range("I134").Value = 5086.3477
$ws.Range("K134").Value = 15259.0431
$ws.Range("M134").Value = -12724.0431

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 298737.34
$ws.Range("I99").Value = 4392.8887
$ws.Range("J99").Value = 629874.9
$ws.Range("K99").Value = 4392.8887
$ws.Range("L99").Value = 629874.9
$ws.Range("M99").Value = -2894.8887
$ws.Range("N99").Value = -632870.9

$ws.Range("H107").Value = 603.0833
$ws.Range("I107").Value = 373.5
$ws.Range("K107").Value = 373.5
$ws.Range("M107").Value = 1546.5

$ws.Range("H126").Value = 298737.34
$ws.Range("I126").Value = 4392.8887
$ws.Range("J126").Value = 629874.9
$ws.Range("K126").Value = 13178.6661
$ws.Range("L126").Value = 1889624.7
$ws.Range("M126").Value = -10708.6661
$ws.Range("N126").Value = -1894564.7

$ws.Range("H134").Value = 280211.22
$ws.Range("I134").Value = 2550.1177
$ws.Range("K134").Value = 7650.353099999999
$ws.Range("M134").Value = -5115.353099999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 385768.2
$ws.Range("I5").Value = 1064.2727
$ws.Range("K5").Value = 3192.8181
$ws.Range("M5").Value = -3080.8181

$ws.Range("H56").Value = 7051.9
$ws.Range("I56").Value = 7051.9
$ws.Range("K56").Value = 7051.9
$ws.Range("M56").Value = -6521.9

$ws.Range("H109").Value = 125652.5
$ws.Range("I109").Value = 745.7143
$ws.Range("J109").Value = 1000000
$ws.Range("K109").Value = 2237.1429
$ws.Range("L109").Value = 3000000
$ws.Range("M109").Value = -1197.1429
$ws.Range("N109").Value = -3002080

$ws.Range("H132").Value = 920664.4399999999
$ws.Range("I132").Value = 167661.67
$ws.Range("J132").Value = 1673667.1
$ws.Range("K132").Value = 1508955.03
$ws.Range("L132").Value = 15063003.9
$ws.Range("M132").Value = -1506425.03
$ws.Range("N132").Value = -15068063.9

$ws.Range("H135").Value = 385768.2
$ws.Range("I135").Value = 1064.2727
$ws.Range("K135").Value = 9578.454299999999
$ws.Range("M135").Value = -7043.454299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8767.615
$ws.Range("I70").Value = 7897.5
$ws.Range("J70").Value = 11668
$ws.Range("K70").Value = 7897.5
$ws.Range("L70").Value = 11668
$ws.Range("M70").Value = -7627.5
$ws.Range("N70").Value = -12208

$ws.Range("H73").Value = 8767.615
$ws.Range("I73").Value = 7897.5
$ws.Range("J73").Value = 11668
$ws.Range("K73").Value = 7897.5
$ws.Range("L73").Value = 11668
$ws.Range("M73").Value = -6961.5
$ws.Range("N73").Value = -13540

$ws.Range("H99").Value = 4685
$ws.Range("I99").Value = 4471
$ws.Range("J99").Value = 4899
$ws.Range("K99").Value = 4471
$ws.Range("L99").Value = 4899
$ws.Range("M99").Value = -2225
$ws.Range("N99").Value = -9391

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 333334500
$ws.Range("I16").Value = 500000740
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 500000740
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -500000570
$ws.Range("N16").Value = -2340

$ws.Range("H132").Value = 5193.5
$ws.Range("J132").Value = 6932.727
$ws.Range("L132").Value = 20798.181
$ws.Range("N132").Value = -25858.181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 8140.8
$ws.Range("I74").Value = 4997
$ws.Range("J74").Value = 8624.462
$ws.Range("K74").Value = 4997
$ws.Range("L74").Value = 8624.462
$ws.Range("M74").Value = -4061
$ws.Range("N74").Value = -10496.462

$ws.Range("H77").Value = 8140.8
$ws.Range("I77").Value = 4997
$ws.Range("J77").Value = 8624.462
$ws.Range("K77").Value = 14991
$ws.Range("L77").Value = 25873.386
$ws.Range("M77").Value = -10311
$ws.Range("N77").Value = -35233.386

$ws.Range("H113").Value = 96.888885
$ws.Range("I113").Value = 96.888885
$ws.Range("K113").Value = 290.666655
$ws.Range("M113").Value = 1879.333345

$ws.Range("H122").Value = 58827530
$ws.Range("I122").Value = 90912420
$ws.Range("J122").Value = 5234
$ws.Range("K122").Value = 272737260
$ws.Range("L122").Value = 15702
$ws.Range("M122").Value = -272734810
$ws.Range("N122").Value = -20602

$ws.Range("H136").Value = 11910311
$ws.Range("I136").Value = 15608243
$ws.Range("J136").Value = 288239
$ws.Range("K136").Value = 46824729
$ws.Range("L136").Value = 864717
$ws.Range("M136").Value = -46822179
$ws.Range("N136").Value = -869817
